# Apply crypto price/volume updates per commit "Updated cryptos list on Mon Sep 18 13:36:59 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.390.89"
$ws.Range("E2").Value = "  +2.12%  "

# Row 3
$ws.Range("D3").Value = "1.665.24"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.61%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.29"
$ws.Range("E5").Value = "  +1.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7
$ws.Range("E7").Value = "  -0.63%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +0.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  +3.87%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").Value = "1.898.85"
$ws.Range("E12").Value = "  +1.31%  "

# Row 13
$ws.Range("D13").Value = "1.661.96"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("E14").Value = "  +1.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  +1.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.28"
$ws.Range("E16").Value = "  +4.00%  "

# Row 17
$ws.Range("D17").Value = "27.357.14"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.63"
$ws.Range("E19").Value = "  +4.85%  "

# Row 20
$ws.Range("E20").Value = "  -0.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +8.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("E22").Value = "  +1.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.49"
$ws.Range("E23").Value = "  +4.18%  "

# Row 24
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.79"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26
$ws.Range("E26").Value = "  -0.78%  "

# Row 27
$ws.Range("E27").Value = "  +3.46%  "

# Row 28
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.05"
$ws.Range("E29").Value = "  +2.62%  "

# Row 30
$ws.Range("E30").Value = "  +0.98%  "

# Row 31
$ws.Range("E31").Value = "  +1.50%  "

# Row 32
$ws.Range("E32").Value = "  +1.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +0.57%  "

# Row 34
$ws.Range("E34").Value = "  +2.38%  "

# Row 35
$ws.Range("D35").Value = "1.270.01"
$ws.Range("E35").Value = "  -1.23%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("E37").Value = "  -0.47%  "

# Row 38
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("E39").Value = "  +1.63%  "

# Row 40
$ws.Range("E40").Value = "  -0.59%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("E42").Value = "  +1.68%  "

# Row 43
$ws.Range("D43").Value = "1.810.68"
$ws.Range("E43").Value = "  +1.37%  "

# Row 44
$ws.Range("E44").Value = "  -4.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.86"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.73"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
$ws.Range("E49").Value = "  +1.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("E51").Value = "  +0.05%  "
